$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "11÷7="
$t.Cell(1, 2).Range.Text = "16÷9="
$t.Cell(1, 3).Range.Text = "16÷7="
$t.Cell(1, 4).Range.Text = "56÷7="
$t.Cell(1, 5).Range.Text = "95÷7="
$t.Cell(5, 1).Range.Text = "22÷2="
$t.Cell(5, 2).Range.Text = "10÷6="
$t.Cell(5, 3).Range.Text = "63÷6="
$t.Cell(5, 4).Range.Text = "35÷5="
$t.Cell(5, 5).Range.Text = "84÷8="
$t.Cell(9, 1).Range.Text = "57÷2="
$t.Cell(9, 2).Range.Text = "85÷7="
$t.Cell(9, 3).Range.Text = "95÷4="
$t.Cell(9, 4).Range.Text = "82÷7="
$t.Cell(9, 5).Range.Text = "73÷9="
$t.Cell(13, 1).Range.Text = "71÷4="
$t.Cell(13, 2).Range.Text = "82÷8="
$t.Cell(13, 3).Range.Text = "71÷9="
$t.Cell(13, 4).Range.Text = "73÷7="
$t.Cell(13, 5).Range.Text = "90÷7="
$t.Cell(17, 1).Range.Text = "35÷2="
$t.Cell(17, 2).Range.Text = "48÷3="
$t.Cell(17, 3).Range.Text = "71÷6="
$t.Cell(17, 4).Range.Text = "88÷5="
$t.Cell(17, 5).Range.Text = "50÷6="

Write-Host "Done updating table cells."
